# Auto-generated script applying numeric updates to Gungnir_Profits sheets
# per the scheduled runner profit recalculation.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Sheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1704.4314
$ws.Cells.Item(15, 9).Value = 1704.4314
$ws.Cells.Item(15, 11).Value = 5113.2942
$ws.Cells.Item(15, 13).Value = -4944.2942
$ws.Cells.Item(62, 8).Value = 9832.35
$ws.Cells.Item(62, 9).Value = 9091.786
$ws.Cells.Item(62, 10).Value = 11560.333
$ws.Cells.Item(62, 11).Value = 9091.786
$ws.Cells.Item(62, 12).Value = 11560.333
$ws.Cells.Item(62, 13).Value = -8467.786
$ws.Cells.Item(62, 14).Value = -12808.333
$ws.Cells.Item(64, 8).Value = 3475
$ws.Cells.Item(64, 9).Value = 3328.5715
$ws.Cells.Item(64, 10).Value = 3680
$ws.Cells.Item(64, 11).Value = 3328.5715
$ws.Cells.Item(64, 12).Value = 3680
$ws.Cells.Item(64, 13).Value = -3080.5715
$ws.Cells.Item(64, 14).Value = -4176
$ws.Cells.Item(65, 8).Value = 9832.35
$ws.Cells.Item(65, 9).Value = 9091.786
$ws.Cells.Item(65, 10).Value = 11560.333
$ws.Cells.Item(65, 11).Value = 45458.93
$ws.Cells.Item(65, 12).Value = 57801.665
$ws.Cells.Item(65, 13).Value = -42338.93
$ws.Cells.Item(65, 14).Value = -64041.665
$ws.Cells.Item(67, 8).Value = 3475
$ws.Cells.Item(67, 9).Value = 3328.5715
$ws.Cells.Item(67, 10).Value = 3680
$ws.Cells.Item(67, 11).Value = 3328.5715
$ws.Cells.Item(67, 12).Value = 3680
$ws.Cells.Item(67, 13).Value = -2470.5715
$ws.Cells.Item(67, 14).Value = -5396
$ws.Cells.Item(74, 8).Value = 3542.8667
$ws.Cells.Item(74, 9).Value = 3740.6
$ws.Cells.Item(74, 10).Value = 3444
$ws.Cells.Item(74, 11).Value = 3740.6
$ws.Cells.Item(74, 12).Value = 3444
$ws.Cells.Item(74, 13).Value = -2804.6
$ws.Cells.Item(74, 14).Value = -5316
$ws.Cells.Item(76, 8).Value = 10422791
$ws.Cells.Item(76, 9).Value = 10646.154
$ws.Cells.Item(76, 10).Value = 17546890
$ws.Cells.Item(76, 11).Value = 10646.154
$ws.Cells.Item(76, 12).Value = 17546890
$ws.Cells.Item(76, 13).Value = -10331.154
$ws.Cells.Item(76, 14).Value = -17547520
$ws.Cells.Item(77, 8).Value = 3542.8667
$ws.Cells.Item(77, 9).Value = 3740.6
$ws.Cells.Item(77, 10).Value = 3444
$ws.Cells.Item(77, 11).Value = 18703
$ws.Cells.Item(77, 12).Value = 17220
$ws.Cells.Item(77, 13).Value = -14023
$ws.Cells.Item(77, 14).Value = -26580
$ws.Cells.Item(79, 8).Value = 10422791
$ws.Cells.Item(79, 9).Value = 10646.154
$ws.Cells.Item(79, 10).Value = 17546890
$ws.Cells.Item(79, 11).Value = 10646.154
$ws.Cells.Item(79, 12).Value = 17546890
$ws.Cells.Item(79, 13).Value = -9554.154
$ws.Cells.Item(79, 14).Value = -17549074
$ws.Cells.Item(137, 8).Value = 1114.7297
$ws.Cells.Item(137, 9).Value = 1005.5769
$ws.Cells.Item(137, 11).Value = 3016.7307
$ws.Cells.Item(137, 13).Value = -466.7307000000001

# ----- Sheet: ARM -----
$ws = $wb.Sheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15627234
$ws.Cells.Item(32, 9).Value = 1893.2373
$ws.Cells.Item(32, 11).Value = 1893.2373
$ws.Cells.Item(32, 13).Value = -1606.2373
$ws.Cells.Item(61, 8).Value = 7938182
$ws.Cells.Item(61, 9).Value = 8773358
$ws.Cells.Item(61, 11).Value = 8773358
$ws.Cells.Item(61, 13).Value = -8773146
$ws.Cells.Item(63, 8).Value = 2478.6365
$ws.Cells.Item(66, 8).Value = 2478.6365
$ws.Cells.Item(74, 8).Value = 1646.0416
$ws.Cells.Item(74, 9).Value = 1863.7142
$ws.Cells.Item(74, 10).Value = 1341.3
$ws.Cells.Item(74, 11).Value = 1863.7142
$ws.Cells.Item(74, 12).Value = 1341.3
$ws.Cells.Item(74, 13).Value = -989.7141999999999
$ws.Cells.Item(74, 14).Value = -3089.3
$ws.Cells.Item(77, 8).Value = 1646.0416
$ws.Cells.Item(77, 9).Value = 1863.7142
$ws.Cells.Item(77, 10).Value = 1341.3
$ws.Cells.Item(77, 11).Value = 9318.571
$ws.Cells.Item(77, 12).Value = 6706.5
$ws.Cells.Item(77, 13).Value = -4950.571
$ws.Cells.Item(77, 14).Value = -15442.5
$ws.Cells.Item(97, 8).Value = 705.5714
$ws.Cells.Item(97, 9).Value = 714
$ws.Cells.Item(97, 10).Value = 655
$ws.Cells.Item(97, 11).Value = 714
$ws.Cells.Item(97, 12).Value = 655
$ws.Cells.Item(97, 13).Value = -218
$ws.Cells.Item(97, 14).Value = -1647
$ws.Cells.Item(132, 8).Value = 1313.55
$ws.Cells.Item(132, 9).Value = 990.7857
$ws.Cells.Item(132, 11).Value = 2972.3571
$ws.Cells.Item(132, 13).Value = -442.3571000000002
$ws.Cells.Item(136, 8).Value = 7938182
$ws.Cells.Item(136, 9).Value = 8773358
$ws.Cells.Item(136, 11).Value = 26320074
$ws.Cells.Item(136, 13).Value = -26317524

# ----- Sheet: BSM -----
$ws = $wb.Sheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 800
$ws.Cells.Item(11, 9).Value = 800
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 800
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 14).Value = -660
$ws.Cells.Item(11, 14).ClearContents()  # cell removed in target
$ws.Cells.Item(134, 8).Value = 3368180.8
$ws.Cells.Item(134, 9).Value = 1036
$ws.Cells.Item(134, 10).Value = 15874718
$ws.Cells.Item(134, 11).Value = 3108
$ws.Cells.Item(134, 12).Value = 47624154
$ws.Cells.Item(134, 13).Value = -573
$ws.Cells.Item(134, 14).Value = -47629224

# ----- Sheet: CRP -----
$ws = $wb.Sheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1011237
$ws.Cells.Item(31, 9).Value = 1793040.9
$ws.Cells.Item(31, 10).Value = 1407
$ws.Cells.Item(31, 11).Value = 1793040.9
$ws.Cells.Item(31, 12).Value = 1407
$ws.Cells.Item(31, 13).Value = -1792745.9
$ws.Cells.Item(31, 14).Value = -1997
$ws.Cells.Item(34, 8).Value = 1011237
$ws.Cells.Item(34, 9).Value = 1793040.9
$ws.Cells.Item(34, 10).Value = 1407
$ws.Cells.Item(34, 11).Value = 1793040.9
$ws.Cells.Item(34, 12).Value = 1407
$ws.Cells.Item(34, 13).Value = -1792838.9
$ws.Cells.Item(34, 14).Value = -1811
$ws.Cells.Item(58, 8).Value = 35715030
$ws.Cells.Item(58, 9).Value = 40000716
$ws.Cells.Item(58, 10).Value = 1000
$ws.Cells.Item(58, 11).Value = 40000716
$ws.Cells.Item(58, 12).Value = 1000
$ws.Cells.Item(58, 13).Value = -40000513
$ws.Cells.Item(58, 14).Value = -1406
$ws.Cells.Item(62, 8).Value = 3753.8572
$ws.Cells.Item(62, 9).Value = 2716.1667
$ws.Cells.Item(62, 11).Value = 2716.1667
$ws.Cells.Item(62, 13).Value = -2092.1667
$ws.Cells.Item(65, 8).Value = 3753.8572
$ws.Cells.Item(65, 9).Value = 2716.1667
$ws.Cells.Item(65, 11).Value = 13580.8335
$ws.Cells.Item(65, 13).Value = -10460.8335
$ws.Cells.Item(99, 8).Value = 47621000
$ws.Cells.Item(99, 9).Value = 111112824
$ws.Cells.Item(99, 10).Value = 2137
$ws.Cells.Item(99, 11).Value = 111112824
$ws.Cells.Item(99, 12).Value = 2137
$ws.Cells.Item(99, 13).Value = -111111326
$ws.Cells.Item(99, 14).Value = -5133
$ws.Cells.Item(126, 8).Value = 47621000
$ws.Cells.Item(126, 9).Value = 111112824
$ws.Cells.Item(126, 10).Value = 2137
$ws.Cells.Item(126, 11).Value = 333338472
$ws.Cells.Item(126, 12).Value = 6411
$ws.Cells.Item(126, 13).Value = -333336002
$ws.Cells.Item(126, 14).Value = -11351
$ws.Cells.Item(132, 8).Value = 11906197
$ws.Cells.Item(132, 9).Value = 1087.1177
$ws.Cells.Item(132, 10).Value = 30305004
$ws.Cells.Item(132, 11).Value = 3261.3531
$ws.Cells.Item(132, 12).Value = 90915012
$ws.Cells.Item(132, 13).Value = -731.3531000000003
$ws.Cells.Item(132, 14).Value = -90920072
$ws.Cells.Item(134, 8).Value = 1403.0588
$ws.Cells.Item(134, 9).Value = 1294.1538
$ws.Cells.Item(134, 10).Value = 1757
$ws.Cells.Item(134, 11).Value = 3882.4614
$ws.Cells.Item(134, 12).Value = 5271
$ws.Cells.Item(134, 13).Value = -1347.4614
$ws.Cells.Item(134, 14).Value = -10341
$ws.Cells.Item(136, 8).Value = 35715030
$ws.Cells.Item(136, 9).Value = 40000716
$ws.Cells.Item(136, 10).Value = 1000
$ws.Cells.Item(136, 11).Value = 120002148
$ws.Cells.Item(136, 12).Value = 3000
$ws.Cells.Item(136, 13).Value = -119999598
$ws.Cells.Item(136, 14).Value = -8100

# ----- Sheet: CUL -----
$ws = $wb.Sheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 27782492
$ws.Cells.Item(5, 9).Value = 37037416
$ws.Cells.Item(5, 10).Value = 17714.166
$ws.Cells.Item(5, 11).Value = 111112248
$ws.Cells.Item(5, 12).Value = 53142.49800000001
$ws.Cells.Item(5, 13).Value = -111112136
$ws.Cells.Item(5, 14).Value = -53366.49800000001
$ws.Cells.Item(51, 8).Value = 2125.125
$ws.Cells.Item(51, 9).Value = 1000
$ws.Cells.Item(51, 10).Value = 2285.8572
$ws.Cells.Item(51, 11).Value = 3000
$ws.Cells.Item(51, 12).Value = 6857.571599999999
$ws.Cells.Item(51, 13).Value = -2540  # newly added cell
$ws.Cells.Item(51, 14).Value = -7777.571599999999
$ws.Cells.Item(115, 8).Value = 10677
$ws.Cells.Item(115, 9).Value = 443.5
$ws.Cells.Item(115, 10).Value = 17499.334
$ws.Cells.Item(115, 11).Value = 1330.5
$ws.Cells.Item(115, 12).Value = 52498.00199999999
$ws.Cells.Item(115, 13).Value = -155.5
$ws.Cells.Item(115, 14).Value = -54848.00199999999
$ws.Cells.Item(122, 8).Value = 11578424
$ws.Cells.Item(122, 9).Value = 44643140
$ws.Cells.Item(122, 10).Value = 5773.325
$ws.Cells.Item(122, 11).Value = 401788260
$ws.Cells.Item(122, 12).Value = 51959.925
$ws.Cells.Item(122, 13).Value = -401785810
$ws.Cells.Item(122, 14).Value = -56859.925
$ws.Cells.Item(131, 8).Value = 780.65
$ws.Cells.Item(131, 9).Value = 551.1429000000001
$ws.Cells.Item(131, 10).Value = 797.92474
$ws.Cells.Item(131, 11).Value = 1653.4287
$ws.Cells.Item(131, 12).Value = 2393.77422
$ws.Cells.Item(131, 13).Value = 3386.5713
$ws.Cells.Item(131, 14).Value = -12473.77422
$ws.Cells.Item(135, 8).Value = 27782492
$ws.Cells.Item(135, 9).Value = 37037416
$ws.Cells.Item(135, 10).Value = 17714.166
$ws.Cells.Item(135, 11).Value = 333336744
$ws.Cells.Item(135, 12).Value = 159427.494
$ws.Cells.Item(135, 13).Value = -333334209
$ws.Cells.Item(135, 14).Value = -164497.494

# ----- Sheet: GSM -----
$ws = $wb.Sheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 15523
$ws.Cells.Item(70, 9).Value = 28752
$ws.Cells.Item(70, 10).Value = 4939.8
$ws.Cells.Item(70, 11).Value = 28752
$ws.Cells.Item(70, 12).Value = 4939.8
$ws.Cells.Item(70, 13).Value = -28482
$ws.Cells.Item(70, 14).Value = -5479.8
$ws.Cells.Item(73, 8).Value = 15523
$ws.Cells.Item(73, 9).Value = 28752
$ws.Cells.Item(73, 10).Value = 4939.8
$ws.Cells.Item(73, 11).Value = 28752
$ws.Cells.Item(73, 12).Value = 4939.8
$ws.Cells.Item(73, 13).Value = -27816
$ws.Cells.Item(73, 14).Value = -6811.8
$ws.Cells.Item(80, 8).Value = 14287886
$ws.Cells.Item(80, 9).Value = 2540
$ws.Cells.Item(80, 10).Value = 50001250
$ws.Cells.Item(80, 11).Value = 2540
$ws.Cells.Item(80, 12).Value = 50001250
$ws.Cells.Item(80, 13).Value = -1542
$ws.Cells.Item(80, 14).Value = -50003246
$ws.Cells.Item(83, 8).Value = 14287886
$ws.Cells.Item(83, 9).Value = 2540
$ws.Cells.Item(83, 10).Value = 50001250
$ws.Cells.Item(83, 11).Value = 12700
$ws.Cells.Item(83, 12).Value = 250006250
$ws.Cells.Item(83, 13).Value = -7708
$ws.Cells.Item(83, 14).Value = -250016234
$ws.Cells.Item(132, 8).Value = 3925.25
$ws.Cells.Item(132, 9).Value = 3468
$ws.Cells.Item(132, 10).Value = 4199.6
$ws.Cells.Item(132, 11).Value = 10404
$ws.Cells.Item(132, 12).Value = 12598.8
$ws.Cells.Item(132, 13).Value = -7874
$ws.Cells.Item(132, 14).Value = -17658.8
$ws.Cells.Item(135, 8).Value = 63375
$ws.Cells.Item(135, 10).Value = 63375
$ws.Cells.Item(135, 12).Value = 63375
$ws.Cells.Item(135, 14).Value = -73515

# ----- Sheet: WVR -----
$ws = $wb.Sheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 23393.191
$ws.Cells.Item(132, 9).Value = 33333.812
$ws.Cells.Item(132, 10).Value = 2186.5334
$ws.Cells.Item(132, 11).Value = 100001.436
$ws.Cells.Item(132, 12).Value = 6559.600199999999
$ws.Cells.Item(132, 13).Value = -97471.43599999999
$ws.Cells.Item(132, 14).Value = -11619.6002
$ws.Cells.Item(136, 8).Value = 12822311
$ws.Cells.Item(136, 9).Value = 19231826
$ws.Cells.Item(136, 10).Value = 3283.077
$ws.Cells.Item(136, 11).Value = 57695478
$ws.Cells.Item(136, 12).Value = 9849.231
$ws.Cells.Item(136, 13).Value = -57692928
$ws.Cells.Item(136, 14).Value = -14949.231

